$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = 44247.55592337526
}
for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = 44247.53459866898
}
for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = 44247.51330931713
}
